# Update cryptocurrency Price (D) and Volume/1h change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.615.12"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.960.81"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'244.39"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'58.64"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.377"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "'22.24"
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").Value = "2.249.10"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "'0.827"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "'13.69"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "1.985.94"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "36.469.56"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'69.86"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "'228.66"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("D26").Value = "'9.21"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'0.139"
$ws.Range("E27").Value = "  +6.52%  "
$ws.Range("D28").Value = "'160.08"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").Value = "'19.44"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'4.72"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").Value = "'0.0619"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "'4.29"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'2.25"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("D37").Value = "'3.36"
$ws.Range("E37").Value = "  +10.55%  "
$ws.Range("D38").Value = "'5.95"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'0.0984"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").Value = "'15.97"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "1.362.24"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'87.75"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'7.14"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "2.139.62"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "'43.67"
$ws.Range("E51").Value = "  -3.93%  "
